$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule row 11's name/label changes from "R40" to "1".
# Force text storage (the label is a string, even though "1" looks numeric)
# so the cell keeps its original text semantics instead of becoming a number.
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
